$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column A. This shifts the existing
# ID/Title/Rating/Review/Year columns one to the right (B..F) and
# preserves their <col> width/bestFit metadata automatically.
$ws.Columns("A").Insert()

# Column A becomes the new "id" (GUID) column; relabel the old "ID"
# header (now in column B) to "MovieId".
$ws.Range("B1").Value = "MovieId"

$ws.Range("A1").Value = "id"
$ws.Range("A2").Value = "d861be4a-de63-49ba-94e0-57486b060d90"
$ws.Range("A3").Value = "bd7d27f8-2f3d-4044-8117-e9e71e351339"
$ws.Range("A4").Value = "a7730d44-e048-4879-b571-a8a92a94c1be"
$ws.Range("A5").Value = "de44443e-7c36-4a51-8101-be42d0b572a1"
$ws.Range("A6").Value = "86ae4a99-30aa-42f3-bf6e-0d08e535ff7d"
$ws.Range("A7").Value = "a264bd90-22d5-47b7-aefe-0f6df48de7a3"
$ws.Range("A8").Value = "572eedce-7e46-4d32-915a-f07c529fed2d"

# Drop the inherited cell styling so the data cells go back to the
# workbook default format (matches a plain retype / paste-values).
$ws.Cells.ClearFormats()

# Autosize the new id column to fit the GUID strings.
$ws.Columns("A").AutoFit()

# Move the active selection down to C9, under the Title column.
$ws.Range("C9").Select()
